# Update 想去人数 (F column) values across sheets 展览, 演出, 全部类型
# per gh-pages data refresh commit 456a3b4.
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 9332
$ws1.Range("F5").Value = 577
$ws1.Range("F8").Value = 254
$ws1.Range("F9").Value = 341
$ws1.Range("F10").Value = 399
$ws1.Range("F11").Value = 150
$ws1.Range("F12").Value = 164
$ws1.Range("F14").Value = 428
$ws1.Range("F15").Value = 11944
$ws1.Range("F21").Value = 231
$ws1.Range("F23").Value = 114
$ws1.Range("F24").Value = 152
$ws1.Range("F25").Value = 2714
$ws1.Range("F27").Value = 61
$ws1.Range("F28").Value = 6
$ws1.Range("F29").Value = 52
$ws1.Range("F30").Value = 2141
$ws1.Range("F31").Value = 987
$ws1.Range("F32").Value = 4180
$ws1.Range("F33").Value = 3610
$ws1.Range("F34").Value = 459
$ws1.Range("F35").Value = 2617
$ws1.Range("F37").Value = 12
$ws1.Range("F38").Value = 1309
$ws1.Range("F39").Value = 189
$ws1.Range("F41").Value = 96
$ws1.Range("F43").Value = 491
$ws1.Range("F44").Value = 61
$ws1.Range("F46").Value = 211
$ws1.Range("F49").Value = 129

# Sheet 2: 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F10").Value = 15

# Sheet 4: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 9332
$ws4.Range("F8").Value = 577
$ws4.Range("F9").Value = 40
$ws4.Range("F11").Value = 254
$ws4.Range("F12").Value = 341
$ws4.Range("F13").Value = 399
$ws4.Range("F14").Value = 164
$ws4.Range("F15").Value = 428
$ws4.Range("F16").Value = 11944
$ws4.Range("F19").Value = 231
$ws4.Range("F23").Value = 114
$ws4.Range("F24").Value = 152
$ws4.Range("F25").Value = 2714
$ws4.Range("F27").Value = 61
$ws4.Range("F28").Value = 6
$ws4.Range("F29").Value = 52
$ws4.Range("F31").Value = 2141
$ws4.Range("F32").Value = 987
$ws4.Range("F33").Value = 4180
$ws4.Range("F34").Value = 3610
$ws4.Range("F35").Value = 460
$ws4.Range("F36").Value = 2617
$ws4.Range("F38").Value = 12
$ws4.Range("F39").Value = 1309
$ws4.Range("F40").Value = 189
$ws4.Range("F41").Value = 770
$ws4.Range("F43").Value = 491
$ws4.Range("F44").Value = 61
$ws4.Range("F46").Value = 211
$ws4.Range("F49").Value = 129

